$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2023-03-22-07:00", 28188.9, 28226.7, 28002.6, 28088.2, 28128.9540757772, 28128.9540757772, 20732.155),
    @("2023-03-22-07:00", 28188.9, 28226.7, 28002.6, 28088.2, 28128.9540757772, 28128.9540757772, 20732.155),
    @("2023-03-22-07:00", 28188.9, 28226.7, 28002.6, 28088.2, 28128.9540757772, 28128.9540757772, 20732.155),
    @("2023-03-22-07:00", 28188.9, 28226.7, 28002.6, 28088.2, 28128.9540757772, 28128.9540757772, 20732.155),
    @("2023-03-22-08:00", 28088.2, 28213.7, 28036,   28071.7, 28128.9540757772, 28003.57721159537, 18427.464)
)

$startRow = 53
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
